$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value2 = 3790
$ws.Range("I64").Value2 = 3995
$ws.Range("J64").Value2 = 3653.3333
$ws.Range("K64").Value2 = 3995
$ws.Range("L64").Value2 = 3653.3333
$ws.Range("M64").Value2 = -3747
$ws.Range("N64").Value2 = -4149.3333
$ws.Range("H67").Value2 = 3790
$ws.Range("I67").Value2 = 3995
$ws.Range("J67").Value2 = 3653.3333
$ws.Range("K67").Value2 = 3995
$ws.Range("L67").Value2 = 3653.3333
$ws.Range("M67").Value2 = -3137
$ws.Range("N67").Value2 = -5369.3333
$ws.Range("H76").Value2 = 4656
$ws.Range("I76").Value2 = 3812.25
$ws.Range("J76").Value2 = 5499.75
$ws.Range("K76").Value2 = 3812.25
$ws.Range("L76").Value2 = 5499.75
$ws.Range("M76").Value2 = -3497.25
$ws.Range("N76").Value2 = -6129.75
$ws.Range("H79").Value2 = 4656
$ws.Range("I79").Value2 = 3812.25
$ws.Range("J79").Value2 = 5499.75
$ws.Range("K79").Value2 = 3812.25
$ws.Range("L79").Value2 = 5499.75
$ws.Range("M79").Value2 = -2720.25
$ws.Range("N79").Value2 = -7683.75
$ws.Range("H107").Value2 = 2156.6072
$ws.Range("I107").Value2 = 1574.5
$ws.Range("J107").Value2 = 5649.25
$ws.Range("K107").Value2 = 1574.5
$ws.Range("L107").Value2 = 5649.25
$ws.Range("M107").Value2 = 345.5
$ws.Range("N107").Value2 = -9489.25
$ws.Range("H138").Value2 = 1940.15
$ws.Range("I138").Value2 = 767.4211
$ws.Range("J138").Value2 = 2215.2346
$ws.Range("K138").Value2 = 2302.2633
$ws.Range("L138").Value2 = 6645.703799999999
$ws.Range("M138").Value2 = 2837.7367
$ws.Range("N138").Value2 = -16925.7038

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value2 = 1644.2222
$ws.Range("I88").Value2 = 1358
$ws.Range("K88").Value2 = 1358
$ws.Range("M88").Value2 = -952
$ws.Range("H91").Value2 = 1644.2222
$ws.Range("I91").Value2 = 1358
$ws.Range("K91").Value2 = 1358
$ws.Range("M91").Value2 = 46

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value2 = 250025970
$ws.Range("I105").Value2 = 250025970
$ws.Range("K105").Value2 = 250025970
$ws.Range("M105").Value2 = -250024223

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value2 = 2085
$ws.Range("J41").Value2 = 0
$ws.Range("L41").Value2 = 0
$ws.Range("N41").ClearContents()
$ws.Range("H58").Value2 = 3891.2444
$ws.Range("I58").Value2 = 1105.9615
$ws.Range("J58").Value2 = 7702.684
$ws.Range("K58").Value2 = 1105.9615
$ws.Range("L58").Value2 = 7702.684
$ws.Range("M58").Value2 = -902.9614999999999
$ws.Range("N58").Value2 = -8108.684
$ws.Range("H62").Value2 = 22224322
$ws.Range("J62").Value2 = 40001860
$ws.Range("L62").Value2 = 40001860
$ws.Range("N62").Value2 = -40003108
$ws.Range("H65").Value2 = 22224322
$ws.Range("J65").Value2 = 40001860
$ws.Range("L65").Value2 = 200009300
$ws.Range("N65").Value2 = -200015540
$ws.Range("H122").Value2 = 1292
$ws.Range("I122").Value2 = 1156
$ws.Range("K122").Value2 = 3468
$ws.Range("M122").Value2 = -1018
$ws.Range("H130").Value2 = 38750
$ws.Range("J130").Value2 = 38750
$ws.Range("L130").Value2 = 38750
$ws.Range("N130").Value2 = -48790
$ws.Range("H136").Value2 = 3891.2444
$ws.Range("I136").Value2 = 1105.9615
$ws.Range("J136").Value2 = 7702.684
$ws.Range("K136").Value2 = 3317.8845
$ws.Range("L136").Value2 = 23108.052
$ws.Range("M136").Value2 = -767.8844999999997
$ws.Range("N136").Value2 = -28208.052

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value2 = 460.54053
$ws.Range("I5").Value2 = 288.7586
$ws.Range("J5").Value2 = 1083.25
$ws.Range("K5").Value2 = 866.2758
$ws.Range("L5").Value2 = 3249.75
$ws.Range("M5").Value2 = -754.2758
$ws.Range("N5").Value2 = -3473.75
$ws.Range("H98").Value2 = 632.8
$ws.Range("I98").Value2 = 180.5
$ws.Range("J98").Value2 = 934.3333
$ws.Range("K98").Value2 = 541.5
$ws.Range("L98").Value2 = 2802.9999
$ws.Range("M98").Value2 = 956.5
$ws.Range("N98").Value2 = -5798.9999
$ws.Range("H107").Value2 = 5474.846
$ws.Range("J107").Value2 = 6107.304
$ws.Range("L107").Value2 = 18321.912
$ws.Range("N107").Value2 = -22161.912
$ws.Range("H113").Value2 = 717.38464
$ws.Range("J113").Value2 = 737.04346
$ws.Range("L113").Value2 = 2211.13038
$ws.Range("N113").Value2 = -6551.130380000001
$ws.Range("H117").Value2 = 866.1818
$ws.Range("I117").Value2 = 607.2222
$ws.Range("J117").Value2 = 2031.5
$ws.Range("K117").Value2 = 1821.6666
$ws.Range("L117").Value2 = 6094.5
$ws.Range("M117").Value2 = 1620.3334
$ws.Range("N117").Value2 = -12978.5
$ws.Range("H131").Value2 = 27820328
$ws.Range("I131").Value2 = 66667184
$ws.Range("J131").Value2 = 72572.234
$ws.Range("K131").Value2 = 200001552
$ws.Range("L131").Value2 = 217716.702
$ws.Range("M131").Value2 = -199996512
$ws.Range("N131").Value2 = -227796.702
$ws.Range("H135").Value2 = 460.54053
$ws.Range("I135").Value2 = 288.7586
$ws.Range("J135").Value2 = 1083.25
$ws.Range("K135").Value2 = 2598.8274
$ws.Range("L135").Value2 = 9749.25
$ws.Range("M135").Value2 = -63.82740000000013
$ws.Range("N135").Value2 = -14819.25

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value2 = 23000
$ws.Range("J52").Value2 = 23000
$ws.Range("L52").Value2 = 23000
$ws.Range("N52").Value2 = -23518
$ws.Range("I53").Value2 = 5000
$ws.Range("J53").Value2 = 0
$ws.Range("K53").Value2 = 5000
$ws.Range("L53").Value2 = 0
$ws.Range("M53").Value2 = -4369
$ws.Range("N53").ClearContents()
$ws.Range("H70").Value2 = 25004116
$ws.Range("I70").Value2 = 17861468
$ws.Range("J70").Value2 = 50003376
$ws.Range("K70").Value2 = 17861468
$ws.Range("L70").Value2 = 50003376
$ws.Range("M70").Value2 = -17861198
$ws.Range("N70").Value2 = -50003916
$ws.Range("H73").Value2 = 25004116
$ws.Range("I73").Value2 = 17861468
$ws.Range("J73").Value2 = 50003376
$ws.Range("K73").Value2 = 17861468
$ws.Range("L73").Value2 = 50003376
$ws.Range("M73").Value2 = -17860532
$ws.Range("N73").Value2 = -50005248
$ws.Range("H80").Value2 = 4683.3335
$ws.Range("I80").Value2 = 0
$ws.Range("J80").Value2 = 4683.3335
$ws.Range("K80").Value2 = 0
$ws.Range("L80").Value2 = 4683.3335
$ws.Range("N80").Value2 = -6679.3335
$ws.Range("M80").ClearContents()
$ws.Range("H83").Value2 = 4683.3335
$ws.Range("I83").Value2 = 0
$ws.Range("J83").Value2 = 4683.3335
$ws.Range("K83").Value2 = 0
$ws.Range("L83").Value2 = 23416.6675
$ws.Range("N83").Value2 = -33400.6675
$ws.Range("M83").ClearContents()
$ws.Range("H102").Value2 = 1833.775
$ws.Range("I102").Value2 = 1402.0769
$ws.Range("K102").Value2 = 1402.0769
$ws.Range("M102").Value2 = 219.9231
$ws.Range("H123").Value2 = 21707.455
$ws.Range("J123").Value2 = 21707.455
$ws.Range("L123").Value2 = 21707.455
$ws.Range("N123").Value2 = -26607.455
$ws.Range("H132").Value2 = 2695.442
$ws.Range("I132").Value2 = 2527.5
$ws.Range("J132").Value2 = 3184
$ws.Range("K132").Value2 = 7582.5
$ws.Range("L132").Value2 = 9552
$ws.Range("M132").Value2 = -5052.5
$ws.Range("N132").Value2 = -14612

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value2 = 1152.3334
$ws.Range("I68").Value2 = 1241.6
$ws.Range("J68").Value2 = 1040.75
$ws.Range("K68").Value2 = 1241.6
$ws.Range("L68").Value2 = 1040.75
$ws.Range("M68").Value2 = -492.5999999999999
$ws.Range("N68").Value2 = -2538.75
$ws.Range("H71").Value2 = 1152.3334
$ws.Range("I71").Value2 = 1241.6
$ws.Range("J71").Value2 = 1040.75
$ws.Range("K71").Value2 = 6208
$ws.Range("L71").Value2 = 5203.75
$ws.Range("M71").Value2 = -2464
$ws.Range("N71").Value2 = -12691.75
$ws.Range("H100").Value2 = 1867.4
$ws.Range("I100").Value2 = 1779
$ws.Range("K100").Value2 = 1779
$ws.Range("M100").Value2 = -1238
